#
# Bongani's "Backend" slide: insert a new slide for the backend part of the
# project, right before the existing "Amazon Echo" slide.
#
# The new slide is produced by duplicating the "Amazon Echo" slide (so it
# inherits the same title/content-placeholder layout) and moving the
# duplicate in front of the original; the duplicate is then retitled to
# "Backend" with an empty body, while the original "Amazon Echo" slide is
# left completely untouched (just shifted one position later).
#

$p = $ppt.ActivePresentation

function Get-TitleShape($slide) {
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.PlaceholderFormat.Type -eq 1) {
            return $shape
        }
    }
    return $slide.Shapes.Item(1)
}

# Locate the existing "Amazon Echo" slide.
$echoIndex = 0
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    $titleShape = Get-TitleShape $candidate
    if ($titleShape.TextFrame.TextRange.Text -eq "Amazon Echo") {
        $echoIndex = $i
        break
    }
}

$echoSlide = $p.Slides.Item($echoIndex)

# Duplicate it and move the new copy immediately before the original.
$newSlide = $echoSlide.Duplicate()
$newSlide.MoveTo($echoIndex)

# Retitle the new (now earlier) slide; leave its content placeholder blank,
# and leave the original "Amazon Echo" slide (now one position later)
# completely unchanged.
$backendSlide = $p.Slides.Item($echoIndex)
$backendTitle = Get-TitleShape $backendSlide
$backendTitle.TextFrame.TextRange.Text = "Backend"
